# Q3 Update - 2025
# Applies the refreshed UNHCR UN-RWA data pull:
#  1. Deletes the trailing duplicate row (old row 303, Yemen/2024) so the
#     sheet shrinks from A1:V303 to A1:V302.
#  2. Updates the shared "short-url" value in column B (every data row)
#     from UTv06W to DX3f88.
#  3. Refreshes the updated statistics cells (refugees / asylum_seekers /
#     returned_refugees / stateless / ooc / hst) for rows 289-301.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the last data row (row 303) entirely - its content is an exact
#    duplicate leftover and the refreshed export no longer contains it.
$ws.Rows.Item(303).Delete()

# 2) Update the short-url column (B) for every remaining data row (2-302).
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value2 = "DX3f88"
}

# 3) Refresh the statistics values that changed in the new data pull.
$statUpdates = @{
    289 = @{ N = 61 }
    290 = @{ N = 49597; O = 32; P = 125; T = 2120 }
    291 = @{ N = 6 }
    293 = @{ N = 63977; O = 13899; T = 4215 }
    294 = @{ N = 205; O = 87 }
    295 = @{ N = 73; O = 32 }
    296 = @{ N = 10 }
    297 = @{ T = 765; V = 12134 }
    298 = @{ N = 31; O = 29 }
    299 = @{ N = 38; O = 27 }
    300 = @{ S = 14500 }
    301 = @{ N = 364; O = 80 }
}

foreach ($row in $statUpdates.Keys) {
    $cols = $statUpdates[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value2 = $cols[$col]
    }
}
